# DDAf_2022_Tableau_annexe_Tab16.xlsx -- "Add files via upload" refresh
#
# The underlying OECD source data moved its trailing measurement window
# forward one year (2011-20 -> 2012-21); that changes the six indicator
# column headers in row 2, refreshes most countries'/aggregates' figures
# to the newer vintage, and drops Djibouti's row (now unavailable, shown
# as the sheet's existing ".." placeholder) for this release.
#
# (Note: the source diff also shows the saved-window chrome size in
# xl/workbook.xml's <workbookView> shrinking from 28800x12490 to
# 19200x10400 -- a byproduct of the author's screen/window size at save
# time, not something the Excel object model exposes as a settable,
# persisted property here, so it's left alone.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column header labels: "...2011-20" -> "...2012-21" (row 2, C:H) ---
$ws.Range("C2").Value = "Échelle de satisfaction existentielle Cantril, 2012-21"
$ws.Range("D2").Value = "Affect négatif, mesures les plus récentes 2012-21"
$ws.Range("E2").Value = "Affect positif, mesures les plus récentes 2012-21"
$ws.Range("F2").Value = "Liberté de faire des choix de vie, mesures les plus récentes 2012-21"
$ws.Range("G2").Value = "Générosité, mesures les plus récentes 2012-21"
$ws.Range("H2").Value = "Soutien sociale, mesures les plus récentes 2012-21"

# --- Row 25 (Djibouti): figures withdrawn -> ".." (not-available marker) ---
$ws.Range("C25").Value = ".."
$ws.Range("D25").Value = ".."
$ws.Range("E25").Value = ".."
$ws.Range("F25").Value = ".."
$ws.Range("G25").Value = ".."
$ws.Range("H25").Value = ".."

# --- Refreshed figures for the 2012-21 vintage ---

# row 38 - Afrique de l'Est
$ws.Range("C38").Value = [double]"4.3069809350100403"
$ws.Range("D38").Value = [double]"0.31398930874738001"
$ws.Range("E38").Value = [double]"0.70614951307122997"
$ws.Range("F38").Value = [double]"0.68648891557347003"
$ws.Range("G38").Value = [double]"0.10217988242706"
$ws.Range("H38").Value = [double]"0.70123958858576996"

# row 62 - Afrique
$ws.Range("C62").Value = [double]"4.5010495896035101"
$ws.Range("D62").Value = [double]"0.34835901127217"
$ws.Range("E62").Value = [double]"0.67362115738240003"
$ws.Range("F62").Value = [double]"0.69406166419070003"
$ws.Range("G62").Value = [double]"7.7786147490999996E-3"
$ws.Range("H62").Value = [double]"0.68738903327191003"

# row 63 - Reste du monde
$ws.Range("C63").Value = [double]"5.8937654735845202"
$ws.Range("D63").Value = [double]"0.29445539704627"
$ws.Range("F63").Value = [double]"0.82816123142154996"
$ws.Range("G63").Value = [double]"-4.7972306924999996E-3"

# row 66 - Monde
$ws.Range("C66").Value = [double]"5.47416517520562"
$ws.Range("D66").Value = [double]"0.31080036394058003"
$ws.Range("E66").Value = [double]"0.70677966129394998"
$ws.Range("F66").Value = [double]"0.78775943873020005"
$ws.Range("G66").Value = [double]"-1.1224706609000001E-3"
$ws.Range("H66").Value = [double]"0.80146210678876995"

# row 67 - COMESA
$ws.Range("C67").Value = [double]"4.4254035022523697"
$ws.Range("D67").Value = [double]"0.33511081337929"
$ws.Range("E67").Value = [double]"0.68409664101070999"
$ws.Range("F67").Value = [double]"0.69834480186303005"
$ws.Range("G67").Value = [double]"8.3654248703500002E-3"
$ws.Range("H67").Value = [double]"0.69900141656398995"

# row 68 - CEN-SAD
$ws.Range("C68").Value = [double]"4.6771178245544398"
$ws.Range("D68").Value = [double]"0.36442494595593"
$ws.Range("E68").Value = [double]"0.67578378319739996"
$ws.Range("F68").Value = [double]"0.70023567703637002"
$ws.Range("G68").Value = [double]"1.269218447574E-2"
$ws.Range("H68").Value = [double]"0.65311868895184"

# row 72 - IGAD
$ws.Range("C72").Value = [double]"4.2266583045323696"
$ws.Range("D72").Value = [double]"0.33109540492296002"
$ws.Range("E72").Value = [double]"0.68660712242125999"
$ws.Range("F72").Value = [double]"0.65360676248867999"
$ws.Range("G72").Value = [double]"0.13304734043777"
$ws.Range("H72").Value = [double]"0.70986185471217"

# row 81 - RDM, pays riches en ressources
$ws.Range("C81").Value = [double]"5.6594903048346996"
$ws.Range("D81").Value = [double]"0.28630470440667999"
$ws.Range("F81").Value = [double]"0.81030523075776995"
$ws.Range("G81").Value = [double]"3.4878488410900002E-3"

# row 82 - Afrique (pays riches en ressources exclus)
$ws.Range("C82").Value = [double]"4.4524682118342502"
$ws.Range("D82").Value = [double]"0.34706117403813003"
$ws.Range("E82").Value = [double]"0.67761577704013998"
$ws.Range("F82").Value = [double]"0.71278704282564997"
$ws.Range("G82").Value = [double]"2.1118091822070001E-2"
$ws.Range("H82").Value = [double]"0.67582122561258995"

# row 86 - Afrique, pays a revenu intermediaire, tranche inferieure
$ws.Range("C86").Value = [double]"4.5986801385879499"
$ws.Range("D86").Value = [double]"0.32320380881429001"
$ws.Range("E86").Value = [double]"0.67545933127403002"
$ws.Range("F86").Value = [double]"0.68152459114789998"
$ws.Range("G86").Value = [double]"-1.16575255292E-2"
$ws.Range("H86").Value = [double]"0.69564319849014"

# row 90 - Pays a revenu eleve
$ws.Range("C90").Value = [double]"6.6837128003438302"
$ws.Range("D90").Value = [double]"0.25683281852139001"
$ws.Range("F90").Value = [double]"0.86471313767962998"
$ws.Range("G90").Value = [double]"-2.5535291112500001E-2"

# row 91 - Afrique, pays les moins avances
$ws.Range("C91").Value = [double]"4.3115492853625099"
$ws.Range("D91").Value = [double]"0.36660103499889002"
$ws.Range("E91").Value = [double]"0.66614309056052001"
$ws.Range("F91").Value = [double]"0.68185059880387999"
$ws.Range("G91").Value = [double]"4.3668594797729998E-2"
$ws.Range("H91").Value = [double]"0.66145415244431005"

# row 97 - Afrique, Etats fragiles
$ws.Range("C97").Value = [double]"4.4438498959396799"
$ws.Range("D97").Value = [double]"0.36362232448477"
$ws.Range("E97").Value = [double]"0.67520290974414998"
$ws.Range("F97").Value = [double]"0.67108643687132996"
$ws.Range("G97").Value = [double]"4.2571862423490002E-2"
$ws.Range("H97").Value = [double]"0.68556041609157004"
